$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.915.68"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "1.895.00"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7774"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.52%  "
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3124"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.69"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07363"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08084"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7708"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.492"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.89%  "
$ws.Range("D14").Value = "1.886.06"
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.226"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.25%  "
$ws.Range("D17").Value = "29.902.15"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "247.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007826"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.56%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.148.01"
$ws.Range("E21").Value = "  +0.99%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.106"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1584"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.435"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.025"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.435"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.78%  "
$ws.Range("E31").Value = "  +0.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.481"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05564"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.29%  "
$ws.Range("E34").Value = "  +0.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.238"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7522"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.98%  "
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.683"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01932"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.799"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4469"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "74.06"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.25%  "
$ws.Range("D43").Value = "1.105.00"
$ws.Range("E43").Value = "  +7.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.968"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.86%  "
$ws.Range("E45").Value = "  +1.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9998"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("E47").Value = "  +1.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.50"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.776"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.26%  "
$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.513"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.043"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.36%  "
